# Reorder the "Recorded By" (column G) values so that a leading "System, "
# entry is moved to the end of the comma-separated list instead of the
# beginning, e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # column G
    $val = $cell.Value2

    if ($val -ne $null -and $val -is [string] -and $val.StartsWith("System, ")) {
        $rest = $val.Substring(8)  # remove leading "System, "
        $newVal = $rest + ", System"
        $cell.Value2 = $newVal
    }
}
